# Updates loading_percent values for Case_5_130 (380 kV case) on Sheet1.
# Applies the new simulation results for rows 2-25 (columns B, C, D, E, F, J, L, N).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "B2"; Value = 31.03919159155054 },
    @{ Cell = "C2"; Value = 11.80897108301059 },
    @{ Cell = "D2"; Value = 3.638390504451045 },
    @{ Cell = "E2"; Value = 9.781564457095438 },
    @{ Cell = "F2"; Value = 57.2142934471481 },
    @{ Cell = "J2"; Value = 9.924875537058675 },
    @{ Cell = "L2"; Value = 12.52259699289024 },
    @{ Cell = "N2"; Value = 22.52100569876991 },
    @{ Cell = "B3"; Value = 30.75136037657967 },
    @{ Cell = "C3"; Value = 11.45441418776775 },
    @{ Cell = "D3"; Value = 3.598486616254093 },
    @{ Cell = "E3"; Value = 9.791276608609015 },
    @{ Cell = "F3"; Value = 57.1017890829129 },
    @{ Cell = "J3"; Value = 9.944898989191147 },
    @{ Cell = "L3"; Value = 12.5314107371746 },
    @{ Cell = "N3"; Value = 22.58261734733784 },
    @{ Cell = "B4"; Value = 30.58223800015957 },
    @{ Cell = "C4"; Value = 11.23524848936596 },
    @{ Cell = "D4"; Value = 3.57333032343898 },
    @{ Cell = "E4"; Value = 9.797639032715571 },
    @{ Cell = "F4"; Value = 57.0464265992056 },
    @{ Cell = "J4"; Value = 9.957879266023536 },
    @{ Cell = "L4"; Value = 12.53898630945128 },
    @{ Cell = "N4"; Value = 22.62246927595919 },
    @{ Cell = "B5"; Value = 30.51529977642609 },
    @{ Cell = "C5"; Value = 11.14572470794361 },
    @{ Cell = "D5"; Value = 3.562915946699569 },
    @{ Cell = "E5"; Value = 9.800332367482248 },
    @{ Cell = "F5"; Value = 57.02731752222356 },
    @{ Cell = "J5"; Value = 9.963341812148633 },
    @{ Cell = "L5"; Value = 12.54261745264074 },
    @{ Cell = "N5"; Value = 22.63921799608941 },
    @{ Cell = "B6"; Value = 30.50430624710896 },
    @{ Cell = "C6"; Value = 11.1308512932535 },
    @{ Cell = "D6"; Value = 3.561176855484029 },
    @{ Cell = "E6"; Value = 9.800785676260732 },
    @{ Cell = "F6"; Value = 57.02435295723856 },
    @{ Cell = "J6"; Value = 9.964259326765802 },
    @{ Cell = "L6"; Value = 12.54325325405595 },
    @{ Cell = "N6"; Value = 22.6420298426204 },
    @{ Cell = "B7"; Value = 30.5813271425325 },
    @{ Cell = "C7"; Value = 11.23404178532466 },
    @{ Cell = "D7"; Value = 3.573190529091277 },
    @{ Cell = "E7"; Value = 9.797674948349719 },
    @{ Cell = "F7"; Value = 57.04615491090328 },
    @{ Cell = "J7"; Value = 9.957952234740564 },
    @{ Cell = "L7"; Value = 12.53903307775075 },
    @{ Cell = "N7"; Value = 22.62269309482588 },
    @{ Cell = "B8"; Value = 30.93841035837481 },
    @{ Cell = "C8"; Value = 11.68711554888633 },
    @{ Cell = "D8"; Value = 3.624768659502345 },
    @{ Cell = "E8"; Value = 9.784830507987847 },
    @{ Cell = "F8"; Value = 57.17265778689499 },
    @{ Cell = "J8"; Value = 9.931637605054732 },
    @{ Cell = "L8"; Value = 12.52518683756895 },
    @{ Cell = "N8"; Value = 22.54182985460955 },
    @{ Cell = "B9"; Value = 31.69584672181382 },
    @{ Cell = "C9"; Value = 12.55753325733106 },
    @{ Cell = "D9"; Value = 3.72064164927088 },
    @{ Cell = "E9"; Value = 9.762799133512143 },
    @{ Cell = "F9"; Value = 57.52930836514138 },
    @{ Cell = "J9"; Value = 9.885452127500303 },
    @{ Cell = "L9"; Value = 12.5152040514705 },
    @{ Cell = "N9"; Value = 22.3992895940628 },
    @{ Cell = "B10"; Value = 32.28274083786092 },
    @{ Cell = "C10"; Value = 13.17813083140934 },
    @{ Cell = "D10"; Value = 3.787773898280137 },
    @{ Cell = "E10"; Value = 9.748522915615615 },
    @{ Cell = "F10"; Value = 57.85700943273567 },
    @{ Cell = "J10"; Value = 9.854788709180932 },
    @{ Cell = "L10"; Value = 12.51833087046732 },
    @{ Cell = "N10"; Value = 22.30432079342284 },
    @{ Cell = "B11"; Value = 32.55533601695962 },
    @{ Cell = "C11"; Value = 13.45489887718426 },
    @{ Cell = "D11"; Value = 3.817577881973422 },
    @{ Cell = "E11"; Value = 9.742440155650183 },
    @{ Cell = "F11"; Value = 58.02017144505981 },
    @{ Cell = "J11"; Value = 9.841541836947963 },
    @{ Cell = "L11"; Value = 12.52202057339671 },
    @{ Cell = "N11"; Value = 22.26323205396096 },
    @{ Cell = "B12"; Value = 32.65928184564271 },
    @{ Cell = "C12"; Value = 13.55879991973048 },
    @{ Cell = "D12"; Value = 3.828756638737763 },
    @{ Cell = "E12"; Value = 9.740195734021814 },
    @{ Cell = "F12"; Value = 58.08396225767285 },
    @{ Cell = "J12"; Value = 9.836626004886687 },
    @{ Cell = "L12"; Value = 12.52374311024328 },
    @{ Cell = "N12"; Value = 22.2479765795825 },
    @{ Cell = "B13"; Value = 32.63686462911248 },
    @{ Cell = "C13"; Value = 13.53646482316557 },
    @{ Cell = "D13"; Value = 3.826353894957728 },
    @{ Cell = "E13"; Value = 9.740676489950554 },
    @{ Cell = "F13"; Value = 58.07013496484073 },
    @{ Cell = "J13"; Value = 9.837680256500398 },
    @{ Cell = "L13"; Value = 12.52335767518878 },
    @{ Cell = "N13"; Value = 22.25124859932913 },
    @{ Cell = "B14"; Value = 32.56387377661788 },
    @{ Cell = "C14"; Value = 13.46346559512084 },
    @{ Cell = "D14"; Value = 3.81849972628629 },
    @{ Cell = "E14"; Value = 9.742254324633222 },
    @{ Cell = "F14"; Value = 58.02537950421853 },
    @{ Cell = "J14"; Value = 9.841135397669348 },
    @{ Cell = "L14"; Value = 12.5221557721658 },
    @{ Cell = "N14"; Value = 22.26197088487587 },
    @{ Cell = "B15"; Value = 32.51925592667892 },
    @{ Cell = "C15"; Value = 13.41863058756184 },
    @{ Cell = "D15"; Value = 3.81367478971504 },
    @{ Cell = "E15"; Value = 9.743228470149598 },
    @{ Cell = "F15"; Value = 57.9982259106647 },
    @{ Cell = "J15"; Value = 9.843264841323423 },
    @{ Cell = "L15"; Value = 12.5214619150869 },
    @{ Cell = "N15"; Value = 22.26857817947872 },
    @{ Cell = "B16"; Value = 32.26503162235777 },
    @{ Cell = "C16"; Value = 13.15992236424332 },
    @{ Cell = "D16"; Value = 3.785811204835411 },
    @{ Cell = "E16"; Value = 9.74892870324935 },
    @{ Cell = "F16"; Value = 57.84662818772549 },
    @{ Cell = "J16"; Value = 9.855668509266659 },
    @{ Cell = "L16"; Value = 12.51813530680833 },
    @{ Cell = "N16"; Value = 22.30704856551451 },
    @{ Cell = "B17"; Value = 32.11044926723483 },
    @{ Cell = "C17"; Value = 12.99971306552978 },
    @{ Cell = "D17"; Value = 3.768528366144058 },
    @{ Cell = "E17"; Value = 9.752530882977887 },
    @{ Cell = "F17"; Value = 57.75722243288816 },
    @{ Cell = "J17"; Value = 9.863457227978641 },
    @{ Cell = "L17"; Value = 12.51667479884288 },
    @{ Cell = "N17"; Value = 22.33119015162922 },
    @{ Cell = "B18"; Value = 32.0220713083925 },
    @{ Cell = "C18"; Value = 12.90704874918299 },
    @{ Cell = "D18"; Value = 3.758518668300444 },
    @{ Cell = "E18"; Value = 9.754641511148748 },
    @{ Cell = "F18"; Value = 57.70712620813944 },
    @{ Cell = "J18"; Value = 9.868003204021658 },
    @{ Cell = "L18"; Value = 12.51604821508293 },
    @{ Cell = "N18"; Value = 22.34527463568509 },
    @{ Cell = "B19"; Value = 31.99224236586373 },
    @{ Cell = "C19"; Value = 12.87558900770888 },
    @{ Cell = "D19"; Value = 3.755117747241354 },
    @{ Cell = "E19"; Value = 9.755362794876088 },
    @{ Cell = "F19"; Value = 57.69039309604366 },
    @{ Cell = "J19"; Value = 9.869553762890625 },
    @{ Cell = "L19"; Value = 12.51587274567 },
    @{ Cell = "N19"; Value = 22.35007755613885 },
    @{ Cell = "B20"; Value = 32.12685021492725 },
    @{ Cell = "C20"; Value = 13.01682185374432 },
    @{ Cell = "D20"; Value = 3.770375314414344 },
    @{ Cell = "E20"; Value = 9.752143415777139 },
    @{ Cell = "F20"; Value = 57.76660256960759 },
    @{ Cell = "J20"; Value = 9.862621266867141 },
    @{ Cell = "L20"; Value = 12.51680818452782 },
    @{ Cell = "N20"; Value = 22.32859965452628 },
    @{ Cell = "B21"; Value = 32.58529411942738 },
    @{ Cell = "C21"; Value = 13.48493260495883 },
    @{ Cell = "D21"; Value = 3.82080961310969 },
    @{ Cell = "E21"; Value = 9.741789276966589 },
    @{ Cell = "F21"; Value = 58.03847103884369 },
    @{ Cell = "J21"; Value = 9.840117816055386 },
    @{ Cell = "L21"; Value = 12.52249997742595 },
    @{ Cell = "N21"; Value = 22.25881323997883 },
    @{ Cell = "B22"; Value = 32.88906928533741 },
    @{ Cell = "C22"; Value = 13.78555242829392 },
    @{ Cell = "D22"; Value = 3.853144793930692 },
    @{ Cell = "E22"; Value = 9.735365985914246 },
    @{ Cell = "F22"; Value = 58.22782696162042 },
    @{ Cell = "J22"; Value = 9.825995908763929 },
    @{ Cell = "L22"; Value = 12.52811569771247 },
    @{ Cell = "N22"; Value = 22.21497544384424 },
    @{ Cell = "B23"; Value = 32.72658752428649 },
    @{ Cell = "C23"; Value = 13.62562537331287 },
    @{ Cell = "D23"; Value = 3.83594477387874 },
    @{ Cell = "E23"; Value = 9.738762828522201 },
    @{ Cell = "F23"; Value = 58.1257037129411 },
    @{ Cell = "J23"; Value = 9.833479632038312 },
    @{ Cell = "L23"; Value = 12.52494529305463 },
    @{ Cell = "N23"; Value = 22.23821035487267 },
    @{ Cell = "B24"; Value = 32.11943380272862 },
    @{ Cell = "C24"; Value = 13.00908869990671 },
    @{ Cell = "D24"; Value = 3.769540538106143 },
    @{ Cell = "E24"; Value = 9.752318466067447 },
    @{ Cell = "F24"; Value = 57.76235774552148 },
    @{ Cell = "J24"; Value = 9.862998992614438 },
    @{ Cell = "L24"; Value = 12.51674721710584 },
    @{ Cell = "N24"; Value = 22.32977017923066 },
    @{ Cell = "B25"; Value = 31.48526674141793 },
    @{ Cell = "C25"; Value = 12.32485111582755 },
    @{ Cell = "D25"; Value = 3.69527947371387 },
    @{ Cell = "E25"; Value = 9.768422744724386 },
    @{ Cell = "F25"; Value = 57.42124195434457 },
    @{ Cell = "J25"; Value = 9.897370032262437 },
    @{ Cell = "L25"; Value = 12.51606651179696 },
    @{ Cell = "N25"; Value = 22.43613615576085 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
